# GPE Data Dictionary (metadata) - add two new columns (phone_number, email)
# to the Volunteers table, make the Volunteers table rows visible, hide the
# Entrance_Exame and Exams table rows, and update the autofilter accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Hide the "Entrance_Exame" table rows (39-49) ---
$ws.Range("A39:A49").EntireRow.Hidden = $true

# --- 2. Hide the "Exams" table rows (55-66) ---
$ws.Range("A55:A66").EntireRow.Hidden = $true

# --- 3. Reveal the "Volunteers" table rows (67-74, before the insert) ---
$ws.Range("A67:A74").EntireRow.Hidden = $false

# --- 4. Insert two new rows right after "role_id" (row 72) and before
#        "start_date" (row 73), so the new columns land between them. ---
$ws.Rows.Item(73).Insert()
$ws.Rows.Item(73).Insert()

# --- 5. Fill in the two new Volunteers columns ---
$ws.Range("A73").Value2 = "Volunteers"
$ws.Range("B73").Value2 = "phone_number"
$ws.Range("C73").Value2 = "int"
$ws.Range("D73").Value2 = "Yes"

$ws.Range("A74").Value2 = "Volunteers"
$ws.Range("B74").Value2 = "email"
$ws.Range("C74").Value2 = "varchar(35)"
$ws.Range("D74").Value2 = "Yes"

# --- 6. Update the defined name / autofilter database range to the new
#        extent ($A$1:$F$81, since two rows were added). Clear the existing
#        autofilter first so the range actually grows. ---
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F81").AutoFilter(1, @("Volunteers"), 7)

$filterDbName = $wb.Names.Item(1)
$filterDbName.RefersTo = "=Metadata!`$A`$1:`$F`$81"
